$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 12 datetime value (time-of-day correction)
$ws.Range("A12").Value = 38991.45833333334

# Update existing rows 209, 210, 212 (open/high/low/close values revised)
$ws.Range("C209").Value = 9059276168200
$ws.Range("D209").Value = 9059276168200
$ws.Range("E209").Value = 9059276168200
$ws.Range("F209").Value = 9059276168200

$ws.Range("C210").Value = 9505999258000
$ws.Range("D210").Value = 9505999258000
$ws.Range("E210").Value = 9505999258000
$ws.Range("F210").Value = 9505999258000

$ws.Range("C212").Value = 10834532376000
$ws.Range("D212").Value = 10834532376000
$ws.Range("E212").Value = 10834532376000
$ws.Range("F212").Value = 10834532376000

# Add new row 213 with latest data point
$ws.Range("A213").Value = 45108.41666666666
$ws.Range("A212").Copy()
$ws.Range("A213").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B213").Value = "ECONOMICS:TRM2"

$ws.Range("C213").Value = 11478416920800
$ws.Range("D213").Value = 11478416920800
$ws.Range("E213").Value = 11478416920800
$ws.Range("F213").Value = 11478416920800
$ws.Range("G213").Value = 0
